# Scheduled market-data refresh: update crafting-leve profit figures
# (currentAveragePrice / LevePrice / LeveProfit columns) per sheet with
# the latest pulled values. Generated from the upstream data-refresh diff.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 6410.544
$ws.Range("I15").Value = 6410.544
$ws.Range("K15").Value = 19231.632
$ws.Range("M15").Value = -19062.632
$ws.Range("H17").Value = 340.69388
$ws.Range("J17").Value = 345.70834
$ws.Range("L17").Value = 1037.12502
$ws.Range("N17").Value = -1373.12502
$ws.Range("H62").Value = 2793.0715
$ws.Range("I62").Value = 2325
$ws.Range("J62").Value = 2871.0833
$ws.Range("K62").Value = 2325
$ws.Range("L62").Value = 2871.0833
$ws.Range("M62").Value = -1701
$ws.Range("N62").Value = -4119.0833
$ws.Range("H65").Value = 2793.0715
$ws.Range("I65").Value = 2325
$ws.Range("J65").Value = 2871.0833
$ws.Range("K65").Value = 11625
$ws.Range("L65").Value = 14355.4165
$ws.Range("M65").Value = -8505
$ws.Range("N65").Value = -20595.4165
$ws.Range("H74").Value = 4642.857
$ws.Range("I74").Value = 4700
$ws.Range("K74").Value = 4700
$ws.Range("M74").Value = -3764
$ws.Range("H77").Value = 4642.857
$ws.Range("I77").Value = 4700
$ws.Range("K77").Value = 23500
$ws.Range("M77").Value = -18820
$ws.Range("H80").Value = 777.2727
$ws.Range("I80").Value = 791.25
$ws.Range("J80").Value = 740
$ws.Range("K80").Value = 2373.75
$ws.Range("L80").Value = 2220
$ws.Range("M80").Value = -1375.75
$ws.Range("N80").Value = -4216
$ws.Range("H83").Value = 777.2727
$ws.Range("I83").Value = 791.25
$ws.Range("J83").Value = 740
$ws.Range("K83").Value = 7121.25
$ws.Range("L83").Value = 6660
$ws.Range("M83").Value = -2129.25
$ws.Range("N83").Value = -16644
$ws.Range("H88").Value = 2959.55
$ws.Range("I88").Value = 1331.3334
$ws.Range("J88").Value = 3657.3572
$ws.Range("K88").Value = 1331.3334
$ws.Range("L88").Value = 3657.3572
$ws.Range("M88").Value = -925.3334
$ws.Range("N88").Value = -4469.3572
$ws.Range("H91").Value = 2959.55
$ws.Range("I91").Value = 1331.3334
$ws.Range("J91").Value = 3657.3572
$ws.Range("K91").Value = 1331.3334
$ws.Range("L91").Value = 3657.3572
$ws.Range("M91").Value = 72.66660000000002
$ws.Range("N91").Value = -6465.3572

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 6750
$ws.Range("I63").Value = 3000
$ws.Range("J63").Value = 8000
$ws.Range("K63").Value = 3000
$ws.Range("L63").Value = 8000
$ws.Range("M63").Value = -2314
$ws.Range("N63").Value = -9372
$ws.Range("H66").Value = 6750
$ws.Range("I66").Value = 3000
$ws.Range("J66").Value = 8000
$ws.Range("K66").Value = 15000
$ws.Range("L66").Value = 40000
$ws.Range("M66").Value = -11568
$ws.Range("N66").Value = -46864
$ws.Range("H88").Value = 2848.889
$ws.Range("I88").Value = 2905
$ws.Range("J88").Value = 2400
$ws.Range("K88").Value = 2905
$ws.Range("L88").Value = 2400
$ws.Range("M88").Value = -2499
$ws.Range("N88").Value = -3212
$ws.Range("H91").Value = 2848.889
$ws.Range("I91").Value = 2905
$ws.Range("J91").Value = 2400
$ws.Range("K91").Value = 2905
$ws.Range("L91").Value = 2400
$ws.Range("M91").Value = -1501
$ws.Range("N91").Value = -5208
$ws.Range("H96").Value = 27500
$ws.Range("J96").Value = 27500
$ws.Range("L96").Value = 27500
$ws.Range("N96").Value = -32992
$ws.Range("H97").Value = 578.5714
$ws.Range("I97").Value = 591.6667
$ws.Range("J97").Value = 500
$ws.Range("K97").Value = 591.6667
$ws.Range("L97").Value = 500
$ws.Range("M97").Value = -95.66669999999999
$ws.Range("N97").Value = -1492
$ws.Range("H109").Value = 40792.332
$ws.Range("J109").Value = 40792.332
$ws.Range("L109").Value = 40792.332
$ws.Range("N109").Value = -43566.332
$ws.Range("H111").Value = 52517.6
$ws.Range("J111").Value = 52517.6
$ws.Range("L111").Value = 52517.6
$ws.Range("N111").Value = -60697.6

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1793.7646
$ws.Range("I86").Value = 1582.9333
$ws.Range("J86").Value = 3375
$ws.Range("K86").Value = 1582.9333
$ws.Range("L86").Value = 3375
$ws.Range("M86").Value = -459.9332999999999
$ws.Range("N86").Value = -5621
$ws.Range("H89").Value = 1793.7646
$ws.Range("I89").Value = 1582.9333
$ws.Range("J89").Value = 3375
$ws.Range("K89").Value = 7914.666499999999
$ws.Range("L89").Value = 16875
$ws.Range("M89").Value = -2298.666499999999
$ws.Range("N89").Value = -28107
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H94").Value = 1049.3529
$ws.Range("I94").Value = 1083.6428
$ws.Range("J94").Value = 889.3333
$ws.Range("K94").Value = 1083.6428
$ws.Range("L94").Value = 889.3333
$ws.Range("M94").Value = -632.6428000000001
$ws.Range("N94").Value = -1791.3333
$ws.Range("H107").Value = 217886.95
$ws.Range("I107").Value = 294525.72
$ws.Range("K107").Value = 294525.72
$ws.Range("M107").Value = -292605.72

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4018.4666
$ws.Range("I132").Value = 4288.75
$ws.Range("J132").Value = 3709.5715
$ws.Range("K132").Value = 12866.25
$ws.Range("L132").Value = 11128.7145
$ws.Range("M132").Value = -10336.25
$ws.Range("N132").Value = -16188.7145

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 543.3333
$ws.Range("I15").Value = 40
$ws.Range("J15").Value = 606.25
$ws.Range("K15").Value = 120
$ws.Range("L15").Value = 1818.75
$ws.Range("M15").Value = 20
$ws.Range("N15").Value = -2098.75
$ws.Range("H134").Value = 5328.5186
$ws.Range("I134").Value = 3430.8333
$ws.Range("J134").Value = 6846.6665
$ws.Range("K134").Value = 10292.4999
$ws.Range("L134").Value = 20539.9995
$ws.Range("M134").Value = -5222.499899999999
$ws.Range("N134").Value = -30679.9995

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 431177.66
$ws.Range("J21").Value = 2001329
$ws.Range("L21").Value = 2001329
$ws.Range("N21").Value = -2001675
$ws.Range("H30").Value = 431177.66
$ws.Range("J30").Value = 2001329
$ws.Range("L30").Value = 2001329
$ws.Range("N30").Value = -2001539
$ws.Range("H122").Value = 3675
$ws.Range("I122").Value = 5850
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 17550
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -15100
$ws.Range("N122").Value = -9400

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 342.5
$ws.Range("I22").Value = 290
$ws.Range("J22").Value = 380
$ws.Range("K22").Value = 290
$ws.Range("L22").Value = 380
$ws.Range("M22").Value = 5
$ws.Range("N22").Value = -970
$ws.Range("H27").Value = 342.5
$ws.Range("I27").Value = 290
$ws.Range("J27").Value = 380
$ws.Range("K27").Value = 290
$ws.Range("L27").Value = 380
$ws.Range("M27").Value = -183
$ws.Range("N27").Value = -594
$ws.Range("H46").Value = 1294.0834
$ws.Range("I46").Value = 1353.2222
$ws.Range("J46").Value = 1116.6666
$ws.Range("K46").Value = 1353.2222
$ws.Range("L46").Value = 1116.6666
$ws.Range("M46").Value = -1165.2222
$ws.Range("N46").Value = -1492.6666
$ws.Range("H82").Value = 1920
$ws.Range("I82").Value = 1700
$ws.Range("J82").Value = 2066.6667
$ws.Range("K82").Value = 1700
$ws.Range("L82").Value = 2066.6667
$ws.Range("M82").Value = -1339
$ws.Range("N82").Value = -2788.6667
$ws.Range("H85").Value = 1920
$ws.Range("I85").Value = 1700
$ws.Range("J85").Value = 2066.6667
$ws.Range("K85").Value = 1700
$ws.Range("L85").Value = 2066.6667
$ws.Range("M85").Value = -452
$ws.Range("N85").Value = -4562.6667
$ws.Range("H106").Value = 13342.25
$ws.Range("J106").Value = 13342.25
$ws.Range("L106").Value = 13342.25
$ws.Range("N106").Value = -15866.25
$ws.Range("H132").Value = 7531.591
$ws.Range("I132").Value = 10000.143
$ws.Range("K132").Value = 30000.429
$ws.Range("M132").Value = -27470.429

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2792.7058
$ws.Range("I81").Value = 812
$ws.Range("J81").Value = 4553.3335
$ws.Range("K81").Value = 1624
$ws.Range("L81").Value = 9106.666999999999
$ws.Range("M81").Value = -563
$ws.Range("N81").Value = -11228.667
$ws.Range("H84").Value = 2792.7058
$ws.Range("I84").Value = 812
$ws.Range("J84").Value = 4553.3335
$ws.Range("K84").Value = 8120
$ws.Range("L84").Value = 45533.335
$ws.Range("M84").Value = -2816
$ws.Range("N84").Value = -56141.335
$ws.Range("H104").Value = 21624.143
$ws.Range("J104").Value = 21624.143
$ws.Range("L104").Value = 21624.143
$ws.Range("N104").Value = -28612.143
$ws.Range("H110").Value = 30000
$ws.Range("J110").Value = 30000
$ws.Range("L110").Value = 30000
$ws.Range("N110").Value = -38180
